$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Estadisticos 2P": updated Blancos/Reprobados/Aprobados/Por_Apro and
# Promedio for row 2 (group 4AEM) and the Promedio for row 8 (group 4BLCM),
# after the second-partial rescate grades came in.
# ---------------------------------------------------------------------------
$ws2P = $wb.Worksheets.Item("Estadisticos 2P")

$ws2P.Range("D2").Value = 0
$ws2P.Range("E2").Value = 4
$ws2P.Range("F2").Value = 20
$ws2P.Range("G2").Value = 83.33
$ws2P.Range("H2").Value = 7.5

$ws2P.Range("D8").Value = 0
$ws2P.Range("H8").Value = 8.4

# ---------------------------------------------------------------------------
# Sheet "Rescatables": a new rescatable student (VIVANCO VIVANCO, LUIS AARON)
# is added at the top; the VASQUEZ / XOTLANIHUA rows swap order and VASQUEZ's
# Reprobadas count drops from 2 to 1. The sheet only has 6 data rows, so the
# simplest faithful reproduction is to rewrite A2:G7 with the final values
# directly (avoids Insert()'s format-inheritance side effects).
# ---------------------------------------------------------------------------
$wsR = $wb.Worksheets.Item("Rescatables")

$textoComa = "Reacciones químicas, conservación de la materia en la formación de nuevas substancias."
$textoDosPuntos = "Reacciones químicas: conservación de la materia en la formación de nuevas substancias."

# Row 2 (new): VIVANCO VIVANCO, LUIS AARON
$wsR.Range("A2").Value = 23330051920313
$wsR.Range("B2").Value = "VIVANCO"
$wsR.Range("C2").Value = "VIVANCO"
$wsR.Range("D2").Value = "LUIS AARON"
$wsR.Range("E2").Value = $textoDosPuntos
$wsR.Range("F2").Value = "4APM"
$wsR.Range("G2").Value = 4

# Row 3: OLMOS ORTEGA, ANGEL GABRIEL (was row 2)
$wsR.Range("A3").Value = 23330051920263
$wsR.Range("B3").Value = "OLMOS"
$wsR.Range("C3").Value = "ORTEGA"
$wsR.Range("D3").Value = "ANGEL GABRIEL"
$wsR.Range("E3").Value = $textoComa
$wsR.Range("F3").Value = "4ARHM"
$wsR.Range("G3").Value = 3

# Row 4: SANTIAGO GARCIA, URIEL (was row 3)
$wsR.Range("A4").Value = 23330051920045
$wsR.Range("B4").Value = "SANTIAGO"
$wsR.Range("C4").Value = "GARCIA"
$wsR.Range("D4").Value = "URIEL"
$wsR.Range("E4").Value = $textoComa
$wsR.Range("F4").Value = "4BEM"
$wsR.Range("G4").Value = 3

# Row 5: XOTLANIHUA ZEPAHUA, JUAN ALBERTO (was row 5, now ahead of VASQUEZ)
$wsR.Range("A5").Value = 23330051920025
$wsR.Range("B5").Value = "XOTLANIHUA"
$wsR.Range("C5").Value = "ZEPAHUA"
$wsR.Range("D5").Value = "JUAN ALBERTO"
$wsR.Range("E5").Value = $textoComa
$wsR.Range("F5").Value = "4AEM"
$wsR.Range("G5").Value = 2

# Row 6: VASQUEZ ESPINDOLA, JOSUE YAHIR (was row 4; Reprobadas now 1, not 2)
$wsR.Range("A6").Value = 23330051920023
$wsR.Range("B6").Value = "VASQUEZ"
$wsR.Range("C6").Value = "ESPINDOLA"
$wsR.Range("D6").Value = "JOSUE YAHIR"
$wsR.Range("E6").Value = $textoComa
$wsR.Range("F6").Value = "4AEM"
$wsR.Range("G6").Value = 1

# Row 7: TEMOXTLE GARCIA, HUGO ANTONIO (was row 6)
$wsR.Range("A7").Value = 23330051920311
$wsR.Range("B7").Value = "TEMOXTLE"
$wsR.Range("C7").Value = "GARCIA"
$wsR.Range("D7").Value = "HUGO ANTONIO"
$wsR.Range("E7").Value = $textoDosPuntos
$wsR.Range("F7").Value = "4APM"
$wsR.Range("G7").Value = 1
